$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows below the existing worker-detail table (before the blank
# gap / footer), shifting the footer block down.
$ws.Rows("19:21").Insert()

# Duplicate the three worker-detail rows (16:18) into the new rows (19:21)
# for the new "2509" period - copy values + full formatting (this also
# carries row 18's distinctive bottom-border look onto row 21, which is now
# the last row of the table).
$ws.Range("B16:J18").Copy()
$ws.Range("B19").PasteSpecial()

# The new rows are a clone of 16:18 (period 2508). Re-point them at the new
# period.
$ws.Range("E19").Value = "2509"
$ws.Range("E20").Value = "2509"
$ws.Range("E21").Value = "2509"

# Periodo Mora column now reads centered across the whole table (old + new
# rows).
$ws.Range("E16:E21").HorizontalAlignment = -4108

# Update summary fields: total "Valor Mora" doubled (two periods now) and
# the period count incremented.
$ws.Range("E11").Value = 369360
$ws.Range("F13").Value = 2
